$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "timestamp" column (Z) values for data rows 2-48 to reflect
# the new run timestamps recorded in the commit (re-run of the notebook).
$ws.Range("Z2").Value = "2025-10-17T07:09:32.076206"
$ws.Range("Z3").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z4").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z5").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z6").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z7").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z8").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z9").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z10").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z11").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z12").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z13").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z14").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z15").Value = "2025-10-17T07:09:32.079719"
$ws.Range("Z16").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z17").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z18").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z19").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z20").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z21").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z22").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z23").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z24").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z25").Value = "2025-10-17T07:09:32.142823"
$ws.Range("Z26").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z27").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z28").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z29").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z30").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z31").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z32").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z33").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z34").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z35").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z36").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z37").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z38").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z39").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z40").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z41").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z42").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z43").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z44").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z45").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z46").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z47").Value = "2025-10-17T07:09:32.209467"
$ws.Range("Z48").Value = "2025-10-17T07:09:32.209467"

Write-Output "Updated timestamps for rows 2-48 in column Z"
